$wb = $excel.ActiveWorkbook

# Use the first existing sheet as the formatting template for the header
# row and the row-label cell (bold text, thin border, centered).
$template = $wb.Worksheets.Item(1)

# --- Add sheet: FTNC_Average_Demand102 (appended at the end) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "FTNC_Average_Demand102"

$ws3.Range("B1").Value = "In-vehicle"
$ws3.Range("C1").Value = "At-stop"
$ws3.Range("D1").Value = "Extra"
$ws3.Range("E1").Value = "Tardiness"
$ws3.Range("F1").Value = "Total"
$ws3.Range("A2").Value = "FTNC_Average_Demand_10"
$ws3.Range("B2").Value = 2489.529745747145
$ws3.Range("C2").Value = 12889.67819330659
$ws3.Range("D2").Value = 879.0425212515569
$ws3.Range("E2").Value = 17.60506501503527
$ws3.Range("F2").Value = 16275.8555262108

$template.Range("B1:F1").Copy()
$ws3.Range("B1:F1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$ws3.Range("A2").PasteSpecial(-4122)

# --- Add sheet: FTHC_Average_Demand10 (appended at the end) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet2)
$ws4.Name = "FTHC_Average_Demand10"

$ws4.Range("B1").Value = "In-vehicle"
$ws4.Range("C1").Value = "At-stop"
$ws4.Range("D1").Value = "Extra"
$ws4.Range("E1").Value = "Tardiness"
$ws4.Range("F1").Value = "Total"
$ws4.Range("A2").Value = "FTHC_Average_Demand_10"
$ws4.Range("B2").Value = 2559.442964239976
$ws4.Range("C2").Value = 13023.49297188319
$ws4.Range("D2").Value = 646.1927526245161
$ws4.Range("E2").Value = 8.417222010433864
$ws4.Range("F2").Value = 16237.54591075812

$template.Range("B1:F1").Copy()
$ws4.Range("B1:F1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$ws4.Range("A2").PasteSpecial(-4122)
